# Fill blank "Types 2" (column F) cells with "-" for all data rows, and
# update the active selection to B1, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data occupies rows 2-140 (row 1 is the header row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$dataRange = $ws.Range("F2:F" + $lastRow)

# Find every blank cell in the "Types 2" column and stamp it with "-".
# SpecialCells can return a multi-area range; assigning .Value directly
# on a multi-area range only touches the first area in this runtime, so
# iterate the areas explicitly.
$blanks = $dataRange.SpecialCells(4)
foreach ($area in $blanks.Areas) {
    $area.Value = "-"
}

# Match the workbook's recorded selection change (A1 -> B1).
$ws.Range("B1").Select()
